$d = $word.ActiveDocument
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body><w:p><w:r><w:t xml:space="preserve">Datasets: </w:t></w:r></w:p><w:p><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>not</w:t></w:r><w:r><w:t xml:space="preserve"> used: </w:t></w:r><w:r><w:t>CDC: Provisional_COVID-19_Death_Counts_by_Sex_Age_and_State</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r></w:p><w:p><w:r><w:t>CDC: Indicators_of_Anxiety_or_Depression_Based_on_Reported_Frequency_of_Symptoms_During_Last_7_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Days</w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">also called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Anxiety_data</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>C</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>DC(?):</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>us_states</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r5.InsertXML($xml)
Write-Output "inserted"
